$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header in B1 from "building_no" to "No"
$ws.Range("B1").Value = "No"

# 2. Strip stray spaces from the mazemap URLs in column C (rows 2-20)
for ($r = 2; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $old = $cell.Value2
    if ($old) {
        $new = $old -replace '\s+', ''
        if ($new -ne $old) {
            $cell.Value2 = $new
        }
    }
}

# 3. Update the selected cell in the sheet view to T12
$ws.Range("T12").Select()
